# QA Test finish fix bug_UI7.1 fix bug_bugReport_UserInteract1
#
# Updates the load-testing QA checklist: the last test case (row 9, TCL1.1)
# is marked "failed" and gets a comment describing the storage/pagination
# bug, plus tidies up the alignment of that row.

$xlCenter = -4108

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 ("TCL1.1") -------------------------------------------------

# Result (C9): not run -> failed, and center it vertically too.
$resultCell = $ws.Cells.Item(9, 3)
$resultCell.Value2 = "failed"
$resultCell.VerticalAlignment = $xlCenter

# ID (A9): center it both horizontally and vertically (wrap stays on).
$idCell = $ws.Cells.Item(9, 1)
$idCell.HorizontalAlignment = $xlCenter
$idCell.VerticalAlignment = $xlCenter

# Comment (D9): new cell describing the bug.
$commentCell = $ws.Cells.Item(9, 4)
$commentCell.Value2 = "В одном кладовке сохранено 2300 объектов. Пролистывает максимум 700 объектов. Далее приложение выгружается системой из-за слишком много потребленной памяти"
$commentCell.WrapText = $true

# Grow the row so the wrapped comment is fully visible.
$ws.Rows.Item(9).RowHeight = 68

# --- Selection ----------------------------------------------------------
$ws.Range("C14").Select()
